$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.933.29'
$ws.Range('E2').Value = '  +4.52%  '

$ws.Range('D3').Value = '2.486.86'
$ws.Range('E3').Value = '  +3.86%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '572.84'
$ws.Range('E5').Value = '  +3.50%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '172.61'
$ws.Range('E6').Value = '  +9.29%  '

$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.520'
$ws.Range('E8').Value = '  +3.06%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.181'
$ws.Range('E9').Value = '  +11.22%  '

$ws.Range('D10').Value = '2.484.52'
$ws.Range('E10').Value = '  +3.95%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.161'
$ws.Range('E11').Value = '  -1.02%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.340'
$ws.Range('E12').Value = '  +3.72%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.72'
$ws.Range('E13').Value = '  +0.55%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000184'
$ws.Range('E14').Value = '  +7.03%  '

$ws.Range('D15').Value = '70.846.81'
$ws.Range('E15').Value = '  +4.60%  '

$ws.Range('D16').Value = '2.955.27'
$ws.Range('E16').Value = '  +3.93%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '24.67'
$ws.Range('E17').Value = '  +8.30%  '

$ws.Range('D18').Value = '2.497.77'
$ws.Range('E18').Value = '  +3.99%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.67'
$ws.Range('E19').Value = '  +12.46%  '

$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.07'
$ws.Range('E20').Value = '  +7.79%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '348.75'
$ws.Range('E21').Value = '  +5.78%  '

$ws.Range('B22').Value = 'SuiNetwork'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.08'
$ws.Range('E22').Value = '  +12.09%  '

$ws.Range('B23').Value = 'Polkadot'
$ws.Range('C23').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.93'
$ws.Range('E23').Value = '  +4.32%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.13%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '66.97'
$ws.Range('E25').Value = '  +1.48%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.91'
$ws.Range('E26').Value = '  +7.34%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.72'
$ws.Range('E27').Value = '  +7.42%  '

$ws.Range('D28').Value = '2.595.88'
$ws.Range('E28').Value = '  +2.67%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.10%  '

$ws.Range('D30').Value = '0.0₃0897'
$ws.Range('E30').Value = '  +11.63%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.54'
$ws.Range('E31').Value = '  +7.06%  '

$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '483.13'
$ws.Range('E32').Value = '  +14.90%  '

$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.28'
$ws.Range('E33').Value = '  +13.31%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.19%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.65'
$ws.Range('E35').Value = '  +3.70%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '159.11'
$ws.Range('E36').Value = '  -1.16%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.115'
$ws.Range('E37').Value = '  +9.39%  '

$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.21'
$ws.Range('E38').Value = '  +1.26%  '

$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.53'
$ws.Range('E39').Value = '  +4.93%  '

$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.59'
$ws.Range('E41').Value = '  +7.70%  '

$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.57'
$ws.Range('E42').Value = '  +6.96%  '

$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.308'
$ws.Range('E43').Value = '  +4.57%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '38.06'
$ws.Range('E44').Value = '  +1.88%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.12'
$ws.Range('E45').Value = '  +5.22%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.19'
$ws.Range('E46').Value = '  +12.00%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '142.81'
$ws.Range('E47').Value = '  +11.04%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.46'
$ws.Range('E48').Value = '  +5.04%  '

$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.505'
$ws.Range('E49').Value = '  +6.11%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0730'
$ws.Range('E50').Value = '  +2.80%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.571'
$ws.Range('E51').Value = '  +3.64%  '
